$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18.60700139133239
$ws.Range("C2").Value = 11.49738024166299
$ws.Range("E2").Value = 13.66333712621891
$ws.Range("F2").Value = 53.7468604590909
$ws.Range("G2").Value = 3.679250259252292
$ws.Range("J2").Value = 10.44193988442619
$ws.Range("B3").Value = 18.14624720533144
$ws.Range("C3").Value = 11.04973446999831
$ws.Range("E3").Value = 13.83359769724555
$ws.Range("F3").Value = 52.7664983415047
$ws.Range("G3").Value = 3.684448315380768
$ws.Range("J3").Value = 10.35803330355948
$ws.Range("B4").Value = 17.86583301994208
$ws.Range("C4").Value = 10.77139057126274
$ws.Range("E4").Value = 13.94435589693841
$ws.Range("F4").Value = 52.16772199765322
$ws.Range("G4").Value = 3.687796688292354
$ws.Range("J4").Value = 10.30738995199563
$ws.Range("B5").Value = 17.75239521595361
$ws.Range("C5").Value = 10.65731521418825
$ws.Range("E5").Value = 13.99105022003824
$ws.Range("F5").Value = 51.92477485914958
$ws.Range("G5").Value = 3.689200785836463
$ws.Range("J5").Value = 10.28697899017118
$ws.Range("B6").Value = 17.7336154065259
$ws.Range("C6").Value = 10.63834104557387
$ws.Range("E6").Value = 13.9988978577533
$ws.Range("F6").Value = 51.88450528251916
$ws.Range("G6").Value = 3.68943633307159
$ws.Range("J6").Value = 10.28360362915462
$ws.Range("B7").Value = 17.86429950774309
$ws.Range("C7").Value = 10.7698544169506
$ws.Range("E7").Value = 13.94497932353902
$ws.Range("F7").Value = 52.16444089356481
$ws.Range("G7").Value = 3.687815463812741
$ws.Range("J7").Value = 10.30711375700108
$ws.Range("B8").Value = 18.44773910079916
$ws.Range("C8").Value = 11.34388684260852
$ws.Range("E8").Value = 13.7207499693326
$ws.Range("F8").Value = 53.40833075912536
$ws.Range("G8").Value = 3.681010133433518
$ws.Range("J8").Value = 10.4128311775965
$ws.Range("B9").Value = 19.60281406321936
$ws.Range("C9").Value = 12.43288811486344
$ws.Range("E9").Value = 13.33057399403441
$ws.Range("F9").Value = 55.86131012564098
$ws.Range("G9").Value = 3.668899750425004
$ws.Range("J9").Value = 10.62678689850011
$ws.Range("B10").Value = 20.44688421158001
$ws.Range("C10").Value = 13.20004678579336
$ws.Range("E10").Value = 13.07443290610461
$ws.Range("F10").Value = 57.65684494266746
$ws.Range("G10").Value = 3.660742606553483
$ws.Range("J10").Value = 10.78759809336837
$ws.Range("B11").Value = 20.82770449591857
$ws.Range("C11").Value = 13.54006938768322
$ws.Range("E11").Value = 12.96462042816577
$ws.Range("F11").Value = 58.46917615346469
$ws.Range("G11").Value = 3.657189784988994
$ws.Range("J11").Value = 10.86144039273818
$ws.Range("B12").Value = 20.97129256185237
$ws.Range("C12").Value = 13.667411557324
$ws.Range("E12").Value = 12.9240103998952
$ws.Range("F12").Value = 58.77589920451217
$ws.Range("G12").Value = 3.655866919617391
$ws.Range("J12").Value = 10.88949239433764
$ws.Range("B13").Value = 20.94039841890515
$ws.Range("C13").Value = 13.64005102627863
$ws.Range("E13").Value = 12.93271304126222
$ws.Range("F13").Value = 58.70988426409642
$ws.Range("G13").Value = 3.656150824277191
$ws.Range("J13").Value = 10.8834470658275
$ws.Range("B14").Value = 20.83953091722431
$ws.Range("C14").Value = 13.55057497775
$ws.Range("E14").Value = 12.96125984945754
$ws.Range("F14").Value = 58.49442963766482
$ws.Range("G14").Value = 3.657080502026606
$ws.Range("J14").Value = 10.86374650873709
$ws.Range("B15").Value = 20.77766113068603
$ws.Range("C15").Value = 13.49558026276937
$ws.Range("E15").Value = 12.97887266737515
$ws.Range("F15").Value = 58.36233452759233
$ws.Range("G15").Value = 3.657652882280735
$ws.Range("J15").Value = 10.85169072194244
$ws.Range("B16").Value = 20.4219187228815
$ws.Range("C16").Value = 13.17763381017553
$ws.Range("E16").Value = 13.08174513185587
$ws.Range("F16").Value = 57.60364595646102
$ws.Range("G16").Value = 3.660977952551197
$ws.Range("J16").Value = 10.78278537533358
$ws.Range("B17").Value = 20.20275225031196
$ws.Range("C17").Value = 12.98019178171931
$ws.Range("E17").Value = 13.14657812039167
$ws.Range("F17").Value = 57.1368920236755
$ws.Range("G17").Value = 3.663058079733348
$ws.Range("J17").Value = 10.74068385673364
$ws.Range("B18").Value = 20.0764060799808
$ws.Range("C18").Value = 12.86579134489895
$ws.Range("E18").Value = 13.18449907181098
$ws.Range("F18").Value = 56.86802120745475
$ws.Range("G18").Value = 3.66426938828653
$ws.Range("J18").Value = 10.71653356595687
$ws.Range("B19").Value = 20.03358337099008
$ws.Range("C19").Value = 12.82691776230678
$ws.Range("E19").Value = 13.19744652263729
$ws.Range("F19").Value = 56.77692408639076
$ws.Range("G19").Value = 3.664682077179378
$ws.Range("J19").Value = 10.70836820755269
$ws.Range("B20").Value = 20.22611389194808
$ws.Range("C20").Value = 13.0012974339631
$ws.Range("E20").Value = 13.13961119729944
$ws.Range("F20").Value = 57.18662265022112
$ws.Range("G20").Value = 3.662835108611189
$ws.Range("J20").Value = 10.74515893554397
$ws.Range("B21").Value = 20.86917625237739
$ws.Range("C21").Value = 13.57689565715673
$ws.Range("E21").Value = 12.95284846028731
$ws.Range("F21").Value = 58.55773992739669
$ws.Range("G21").Value = 3.656806823753431
$ws.Range("J21").Value = 10.86953068699895
$ws.Range("B22").Value = 21.285770410145
$ws.Range("C22").Value = 13.9447752113447
$ws.Range("E22").Value = 12.83646786577052
$ws.Range("F22").Value = 59.44855528747073
$ws.Range("G22").Value = 3.652998111434929
$ws.Range("J22").Value = 10.95133207897798
$ws.Range("B23").Value = 21.06381534268773
$ws.Range("C23").Value = 13.74922910969984
$ws.Range("E23").Value = 12.89805943890762
$ws.Range("F23").Value = 58.97367293316385
$ws.Range("G23").Value = 3.655018960132893
$ws.Range("J23").Value = 10.90762900004449
$ws.Range("B24").Value = 20.21555315381589
$ws.Range("C24").Value = 12.9917583257837
$ws.Range("E24").Value = 13.14275892676656
$ws.Range("F24").Value = 57.16414105897228
$ws.Range("G24").Value = 3.662935865823955
$ws.Range("J24").Value = 10.74313558215084
$ws.Range("B25").Value = 19.2903603793392
$ws.Range("C25").Value = 12.14338245331098
$ws.Range("E25").Value = 13.43079322888321
$ws.Range("F25").Value = 55.19781955429924
$ws.Range("G25").Value = 3.672045014788252
$ws.Range("J25").Value = 10.56823775568418
